$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 (Highs-BigM (100,100) / NO_SOLUTION / INFEASIBLE)
$ws.Range("E2").Value = 0.000323266
$ws.Range("F2").Value = 0.02228147
$ws.Range("G2").Value = 0.0005715590879475077

# Row 3 (Highs-BigM (100,100) / NO_SOLUTION / INFEASIBLE_POINT)
$ws.Range("E3").Value = 0.00967684
$ws.Range("F3").Value = 0.018272666
$ws.Range("G3").Value = 0.010776944681034482
